$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.472.56"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "2.979.31"
$ws.Range("E3").Value = "  +2.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.21"
$ws.Range("E5").Value = "  +2.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.45"
$ws.Range("E6").Value = "  +1.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.544"
$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.27"
$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0843"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.457.56"
$ws.Range("E13").Value = "  +2.34%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.57"
$ws.Range("E14").Value = "  +0.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.50"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").Value = "2.990.07"
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.977"
$ws.Range("E17").Value = "  +2.55%  "

$ws.Range("D18").Value = "51.476.89"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.36"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +1.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.15"
$ws.Range("E23").Value = "  +1.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.25"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  +3.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  -1.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("E27").Value = "  +17.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.43"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.99"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.112"
$ws.Range("E31").Value = "  +6.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.85"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.87"
$ws.Range("E33").Value = "  -0.87%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.35"
$ws.Range("E34").Value = "  +1.75%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.09"
$ws.Range("E35").Value = "  -2.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0459"
$ws.Range("E36").Value = "  +7.88%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.42"
$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.60"
$ws.Range("E40").Value = "  -6.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  -0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  +2.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.90"
$ws.Range("E43").Value = "  +5.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.29"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.280"
$ws.Range("E45").Value = "  +19.15%  "

$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  +4.17%  "

$ws.Range("D48").Value = "2.042.47"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.24"
$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0358"
$ws.Range("E50").Value = "  +12.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.18"
$ws.Range("E51").Value = "  +2.93%  "
